$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column I header: "Other found locations"
$ws.Cells.Item(1, 9).Value2 = "Other found locations"

# Row 2 (bsp.2008.0030 record)
$ws.Cells.Item(2, 5).Value2 = '[Atas%Jenny%coreGivesNoEmail%1,   Bandy%Kenneth%coreGivesNoEmail%1,   Bradin%Stuart A.%coreGivesNoEmail%1,   Cadwallender%Bruce A.%coreGivesNoEmail%1,   Cinti%Sandro K.%coreGivesNoEmail%1,   Collins%Curtis D.%coreGivesNoEmail%1,   Goldberg%Janet%coreGivesNoEmail%1,   Holmes%Jennifer G.%coreGivesNoEmail%1,   Kim%Christopher%coreGivesNoEmail%1,   Krupansky%Frank%coreGivesNoEmail%1,   Lozon%Marie M.%coreGivesNoEmail%1,   Rodgers%Phillip E.%coreGivesNoEmail%1,   Shlafer%Jean%coreGivesNoEmail%1,   Wagner%Deborah%coreGivesNoEmail%1,   Wilkerson%William M.%coreGivesNoEmail%1,   Wright%Carrie M.%coreGivesNoEmail%1]'
$ws.Cells.Item(2, 6).Value2 = "not found"
$ws.Cells.Item(2, 7).Value2 = "N/A"
$ws.Cells.Item(2, 9).Value2 = ""

# Row 3 (Impact of Three Influenza Epidemics record)
$ws.Cells.Item(3, 5).Value2 = '[ Douglas M.%Fleming%null%2,  Douglas M.%Fleming%null%0]'
$ws.Cells.Item(3, 9).Value2 = ""

# Row 4 (hpu.2016.0078 record) - only the new column is added
$ws.Cells.Item(4, 9).Value2 = ""

# Row 5 (tmi.12532 record)
$ws.Cells.Item(5, 6).Value2 = "not found"
$ws.Cells.Item(5, 7).Value2 = "N/A"
$ws.Cells.Item(5, 9).Value2 = ""
